$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 48 with quarter label "III-2021" and its data values
$ws.Range("A48").Value = "III-2021"

$rowValues = @(8.4, 4.9, 6.9, 9.6, 8.4, 9.2, 8.7, 7.6, 7.2, 8.4, 8.6, 7.1, 8.6, 3, 3.4, 5.8, 9.2)

$col = 2
foreach ($val in $rowValues) {
    $ws.Cells.Item(48, $col).Value = $val
    $col = $col + 1
}
